$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds/statistics values per row, as captured in the source diff.

# Row 2
$ws.Range("L2").Value = 1.2
$ws.Range("M2").Value = 4.5
$ws.Range("N2").Value = 1.67
$ws.Range("O2").Value = 2.2

# Row 3
$ws.Range("J3").Value = 1.05
$ws.Range("L3").Value = 1.41
$ws.Range("M3").Value = 2.7

# Row 4
$ws.Range("J4").Value = 1.07
$ws.Range("L4").Value = 1.41
$ws.Range("M4").Value = 2.7
$ws.Range("X4").Value = 21
$ws.Range("AE4").Value = 8.5

# Row 5
$ws.Range("J5").Value = 1.05
$ws.Range("L5").Value = 1.41
$ws.Range("M5").Value = 2.7
$ws.Range("U5").Value = 9
$ws.Range("X5").Value = 19
$ws.Range("AA5").Value = 6

# Row 6
$ws.Range("G6").Value = 2.05
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 1.07
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = 1.41
$ws.Range("M6").Value = 2.7
$ws.Range("U6").Value = 8.5
$ws.Range("W6").Value = 17
$ws.Range("AD6").Value = 501
$ws.Range("AF6").Value = 19

# Row 7
$ws.Range("G7").Value = 2.45
$ws.Range("I7").Value = 3.1
$ws.Range("J7").Value = 1.07
$ws.Range("L7").Value = 1.41
$ws.Range("M7").Value = 2.7
$ws.Range("AI7").Value = 29

# Row 8
$ws.Range("J8").Value = 1.03
$ws.Range("L8").Value = 1.22

# Row 9
$ws.Range("J9").Value = 1.05
$ws.Range("L9").Value = 1.41
$ws.Range("M9").Value = 2.7
$ws.Range("N9").Value = 2.25
$ws.Range("O9").Value = 1.62

# Row 10
$ws.Range("G10").Value = 2.35
$ws.Range("I10").Value = 3.2
$ws.Range("P10").Value = 1.54
$ws.Range("AF10").Value = 15

# Row 11
$ws.Range("P11").Value = 1.58
$ws.Range("W11").Value = 19

# Row 12
$ws.Range("J12").Value = 1.07
$ws.Range("K12").Value = 9
$ws.Range("P12").Value = 1.47

# Row 13
$ws.Range("P13").Value = 1.54

# Row 14
$ws.Range("G14").Value = 5.25
$ws.Range("H14").Value = 4.1
$ws.Range("I14").Value = 1.57
$ws.Range("P14").Value = 1.22
$ws.Range("Q14").Value = 4
$ws.Range("R14").Value = 1.5
$ws.Range("S14").Value = 2.5
$ws.Range("X14").Value = 34
$ws.Range("AB14").Value = 12
$ws.Range("AE14").Value = 12
$ws.Range("AJ14").Value = 17

# Row 15
$ws.Range("G15").Value = 2.38
$ws.Range("H15").Value = 3.3
$ws.Range("I15").Value = 3
$ws.Range("N15").Value = 1.93
$ws.Range("O15").Value = 1.93
$ws.Range("P15").Value = 1.36
$ws.Range("Q15").Value = 3
$ws.Range("U15").Value = 12
$ws.Range("AE15").Value = 11
$ws.Range("AI15").Value = 23

# Row 16
$ws.Range("G16").Value = 1.44
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 6
$ws.Range("V16").Value = 8.5
$ws.Range("W16").Value = 11
$ws.Range("AE16").Value = 23
$ws.Range("AG16").Value = 19

# Row 17
$ws.Range("K17").Value = 17

# Row 18
$ws.Range("H18").Value = 8
$ws.Range("K18").Value = 19
$ws.Range("T18").Value = 9
$ws.Range("Y18").Value = 29
$ws.Range("AE18").Value = 41
$ws.Range("AF18").Value = 81
$ws.Range("AH18").Value = 251

# Row 19
$ws.Range("G19").Value = 1.65
$ws.Range("H19").Value = 3.7
$ws.Range("I19").Value = 5.25
$ws.Range("L19").Value = 1.44
$ws.Range("M19").Value = 2.63
$ws.Range("R19").Value = 2.38
$ws.Range("S19").Value = 1.53
$ws.Range("U19").Value = 6.5
$ws.Range("V19").Value = 9
$ws.Range("W19").Value = 12
$ws.Range("Z19").Value = 7.5
$ws.Range("AA19").Value = 7.5
$ws.Range("AC19").Value = 101
$ws.Range("AF19").Value = 26
$ws.Range("AG19").Value = 19
$ws.Range("AH19").Value = 67
$ws.Range("AI19").Value = 51

# Row 21
$ws.Range("G21").Value = 1.22
$ws.Range("I21").Value = 15
$ws.Range("R21").Value = 2.5
$ws.Range("S21").Value = 1.5
$ws.Range("AE21").Value = 29
$ws.Range("AF21").Value = 67
$ws.Range("AG21").Value = 41
$ws.Range("AH21").Value = 201
$ws.Range("AJ21").Value = 101

# Row 22
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = 3.75
$ws.Range("N22").Value = 2.35
$ws.Range("O22").Value = 1.57
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 1.73
$ws.Range("U22").Value = 8.5
$ws.Range("V22").Value = 9.5
$ws.Range("AF22").Value = 17

# Row 25
$ws.Range("G25").Value = 2.8
$ws.Range("I25").Value = 2.35

# Row 36
$ws.Range("J36").Value = 1.06
$ws.Range("K36").Value = 10
$ws.Range("L36").Value = 1.3
$ws.Range("M36").Value = 3.4
$ws.Range("N36").Value = 2.03
$ws.Range("O36").Value = 1.83

# Row 38
$ws.Range("G38").Value = 9.75
$ws.Range("H38").Value = 6.2
$ws.Range("N38").Value = 1.28
$ws.Range("O38").Value = 3.5
$ws.Range("R38").Value = 1.74
$ws.Range("S38").Value = 2.01
$ws.Range("T38").Value = 32
$ws.Range("U38").Value = 70
$ws.Range("Z38").Value = 25
$ws.Range("AA38").Value = 12.5
$ws.Range("AB38").Value = 18.5
$ws.Range("AC38").Value = 55
$ws.Range("AE38").Value = 10.5
$ws.Range("AF38").Value = 7.3
$ws.Range("AG38").Value = 8.5
$ws.Range("AH38").Value = 7.1
$ws.Range("AJ38").Value = 18

# Row 39
$ws.Range("H39").Value = 3.5
$ws.Range("I39").Value = 2.75
$ws.Range("L39").Value = 1.18
$ws.Range("M39").Value = 4.5
$ws.Range("N39").Value = 1.62
$ws.Range("O39").Value = 2.25
$ws.Range("P39").Value = 1.3
$ws.Range("Q39").Value = 3.4
$ws.Range("R39").Value = 1.5
$ws.Range("S39").Value = 2.5
$ws.Range("T39").Value = 12
$ws.Range("Y39").Value = 21
$ws.Range("AA39").Value = 7
$ws.Range("AD39").Value = 101
$ws.Range("AE39").Value = 13

# Row 41
$ws.Range("G41").Value = 2.7
$ws.Range("I41").Value = 2.75
$ws.Range("U41").Value = 13

# Row 43
$ws.Range("H43").Value = 3.9
$ws.Range("I43").Value = 1.7
$ws.Range("Y43").Value = 34
$ws.Range("AA43").Value = 8

# Row 44
$ws.Range("G44").Value = 2
$ws.Range("I44").Value = 3.75
$ws.Range("T44").Value = 5.5
$ws.Range("W44").Value = 17
$ws.Range("AC44").Value = 81
$ws.Range("AE44").Value = 8.5
$ws.Range("AG44").Value = 15
$ws.Range("AI44").Value = 41
